$wb = $excel.ActiveWorkbook

# A temporary sheet is added (and later removed) purely so that the
# internal sheetId counter advances past the id that would otherwise be
# reused for the new "fuel" sheet (matches sheetId="18" from the target).
$tmp = $wb.Worksheets.Add()
$tmpName = $tmp.Name

# Activate the sheet that should immediately follow the new sheet, so the
# new sheet gets inserted right before it (i.e. right after
# "asymmetric_sgen" and right before "ext_grid").
$extGrid = $wb.Worksheets.Item("ext_grid")
$extGrid.Activate()

# Insert the new "fuel" worksheet.
$fuel = $wb.Worksheets.Add()
$fuel.Name = "fuel"

# Header row (bold).
$fuel.Range("B1").Value = "gen_type"
$fuel.Range("C1").Value = "index"
$fuel.Range("D1").Value = "fuel"
$fuel.Range("B1:D1").Font.Bold = $true

# Data rows.
for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $fuel.Cells.Item($r, 1).Value = $i
    $fuel.Cells.Item($r, 2).Value = "sgen"
    $fuel.Cells.Item($r, 3).Value = $i
    $fuel.Cells.Item($r, 4).Value = "solar"
}

$fuel.Range("K27").Select()

# Remove the helper temp sheet (re-fetch by name, since the stale
# reference can shift identity after further sheet insertions).
$wb.Worksheets.Item($tmpName).Delete()
